$d = $word.ActiveDocument

# 1. Author paragraph: "T" -> "Tashfia"
$d.Paragraphs.Item(2).Range.Text = "Tashfia"

# 2. Date paragraph: "7/31/2021" -> "8/7/2021"
$d.Paragraphs.Item(3).Range.Text = "8/7/2021"

# 3. Append two new runs at the end of the last paragraph ("Note that the
#    echo = FALSE parameter was added ... generated the plot."):
#    a run containing a single space, then a run containing
#    "## This is a markdown file"
$r1 = $d.Paragraphs.Last.Range
$r1.Collapse(0)
$r1.InsertAfter(" ")

$r2 = $d.Paragraphs.Last.Range
$r2.Collapse(0)
$r2.InsertAfter("## This is a markdown file")
